# Commit: "Changed text from "uSD" to "MicroSD""
#
# Applies two textual changes:
#  1. The cached "datetimeFigureOut" footer date field (shown on the
#     slide master and on every slide layout) moves from 2019/12/14 to
#     2019/12/18.
#  2. The "Rectangle 16" callout on slide 1 changes from "µSD card" to
#     "MicroSD card".

$p = $ppt.ActivePresentation

$oldDate = "2019/12/14"
$newDate = "2019/12/18"

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# The date placeholder shape lives on the slide master itself ...
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

# ... and is repeated on every slide layout that hangs off the master.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholders $layouts.Item($li).Shapes
}

# Slide-level text fix: micro-sign "SD card" -> "MicroSD card".
$oldLabel = [char]0x00B5 + "SD card"
$newLabel = "MicroSD card"

$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldLabel) {
            $tr.Text = $newLabel
        }
    }
}
